# Going through the dataset, updating.
#
# This adds a new "Other found locations" column (I) that records which
# secondary API/source matched each row (PMC / PMC_elsevier / PMC_Springer),
# refreshes several "Authors" (E) cells whose name-separator formatting
# changed (extra padding from the re-run of the matcher), flips row 10 from
# a CORE match back to "not found" / "N/A", and widens the Authors value on
# row 10 to match the same re-run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column I.
$ws.Range('I1').Value = 'Other found locations'

# Row 2, 11, 12 ([] / not found / N/A rows) also gained a blank "Other found
# locations" cell in column I - nothing was found for them either.
$ws.Range('I2').Value = ''

# Row 3 - PMC7349059: Authors re-spaced, matched via plain PMC lookup.
$ws.Range('E3').Value = '[Marco%Colizzi%NULL%1,    Elena%Sironi%NULL%2,    Elena%Sironi%NULL%0,    Federico%Antonini%NULL%1,    Marco Luigi%Ciceri%NULL%1,    Chiara%Bovo%NULL%2,    Chiara%Bovo%NULL%0,    Leonardo%Zoccante%NULL%1]'
$ws.Range('I3').Value = '_PMC'

# Row 4 - PMC7127630: Authors re-spaced, matched via PMC/Elsevier lookup.
$ws.Range('E4').Value = '[Wen Yan%Jiao%NULL%1,    Lin Na%Wang%NULL%1,    Juan%Liu%NULL%1,    Shuan Feng%Fang%NULL%1,    Fu Yong%Jiao%NULL%1,    Massimo%Pettoello-Mantovani%NULL%1,    Eli%Somekh%NULL%1]'
$ws.Range('I4').Value = '_PMC_elsevier'

# Row 5 - PMC7196181: Authors re-spaced, matched via PMC/Springer lookup.
$ws.Range('E5').Value = '[Shuang-Jiang%Zhou%NULL%0,    Li-Gang%Zhang%NULL%0,    Lei-Lei%Wang%NULL%0,    Zhao-Chang%Guo%NULL%0,    Jing-Qi%Wang%NULL%0,    Jin-Cheng%Chen%NULL%0,    Mei%Liu%NULL%0,    Xi%Chen%NULL%0,    Jing-Xu%Chen%chenjx1110@163.com%0]'
$ws.Range('I5').Value = '_PMC_Springer'

# Row 6 - PMC7256340: Authors re-spaced, matched via PMC/Springer lookup.
$ws.Range('E6').Value = '[Xi%Liu%NULL%1,    Wen-Tao%Luo%NULL%1,    Ying%Li%NULL%0,    Chun-Na%Li%NULL%1,    Zhong-Si%Hong%NULL%1,    Hui-Li%Chen%NULL%1,    Fei%Xiao%NULL%1,    Jin-Yu%Xia%xiajinyu@mail.sysu.edu.cn%1]'
$ws.Range('I6').Value = '_PMC_Springer'

# Row 7 - PMC7205689: Authors re-spaced, matched via PMC/Elsevier lookup.
$ws.Range('E7').Value = '[Benjamin%Oosterhoff%Benjamin.oosterhoff@montana.edu%2,    Cara A.%Palmer%NULL%2,    Jenna%Wilson%NULL%2,    Natalie%Shook%NULL%2]'
$ws.Range('I7').Value = '_PMC_elsevier'

# Row 8 - PMC7293436: Authors re-spaced, matched via PMC/Springer lookup.
$ws.Range('E8').Value = '[İsmail%Seçer%ismailsecer84@gmail.com%2,    Sümeyye%Ulaş%NULL%4,    Sümeyye%Ulaş%NULL%0]'
$ws.Range('I8').Value = '_PMC_Springer'

# Row 9 - PMC7151383: Authors re-spaced, matched via PMC/Elsevier lookup.
$ws.Range('E9').Value = '[Fangyuan%Tian%NULL%1,    Hongxia%Li%lihx@xust.edu.cn%1,    Shuicheng%Tian%tiansc@xust.edu.cn%1,    Jie%Yang%NULL%0,    Jiang%Shao%NULL%1,    Chenning%Tian%NULL%1]'
$ws.Range('I9').Value = '_PMC_elsevier'

# Row 10: no longer resolves via CORE - ID/ID Format revert to "not found"/"N/A",
# and the (CORE-sourced) Authors value is re-spaced like the others.
$ws.Range('E10').Value = '[Abrams%D.%coreGivesNoEmail%1,   Antonovsky%A.%coreGivesNoEmail%1,   Bowlby%J.%coreGivesNoEmail%2,   Bowlby%J.%coreGivesNoEmail%0,   Durkheim%\u00c9.%coreGivesNoEmail%1,   Erikson%E. H.%coreGivesNoEmail%1,   Haidt%J.%coreGivesNoEmail%1,   Janoff-Bulman%R.%coreGivesNoEmail%1,   Jerzy%Trzebi\u0144ski%coreGivesNoEmail%1,   Jolanta%Zuzanna Czarnecka%coreGivesNoEmail%1,   Maciej%Caba\u0144ski%coreGivesNoEmail%1,   Seligman%M. E. P.%coreGivesNoEmail%2,   Seligman%M. E. P.%coreGivesNoEmail%0,   Spielberger%C. D.%coreGivesNoEmail%1,   Trzebi\u0144ski%J.%coreGivesNoEmail%1,   Zi\u0119ba%M.%coreGivesNoEmail%2,   Zi\u0119ba%M.%coreGivesNoEmail%0]'
$ws.Range('F10').Value = 'not found'
$ws.Range('G10').Value = 'N/A'
$ws.Range('I10').Value = ''

$ws.Range('I11').Value = ''
$ws.Range('I12').Value = ''

# Row 13 - PMC7196181 (duplicate of row 5): Authors re-spaced, PMC/Springer lookup.
$ws.Range('E13').Value = '[Shuang-Jiang%Zhou%NULL%0,    Li-Gang%Zhang%NULL%0,    Lei-Lei%Wang%NULL%0,    Zhao-Chang%Guo%NULL%0,    Jing-Qi%Wang%NULL%0,    Jin-Cheng%Chen%NULL%0,    Mei%Liu%NULL%0,    Xi%Chen%NULL%0,    Jing-Xu%Chen%chenjx1110@163.com%0]'
$ws.Range('I13').Value = '_PMC_Springer'
